$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 3).Value = "2023-06-06 12:40:00"
}
